# Applies the diff: swaps match data (columns F:V) between several pairs
# of adjacent rows, and appends 6 new match rows (rows 66-71) at the end
# of the sheet, updating the used-range dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows($r1, $r2) {
    $v1 = $ws.Range("F$r1" + ":V$r1").Value2()
    $v2 = $ws.Range("F$r2" + ":V$r2").Value2()
    $ws.Range("F$r1" + ":V$r1").Value = $v2
    $ws.Range("F$r2" + ":V$r2").Value = $v1
}

# Swap the match details (home/away teams, scores, odds, timestamps, url)
# between these adjacent row pairs while keeping Indice/pais/torneio/
# temporada/data_partida (columns A:E) untouched.
Swap-MatchRows 18 19
Swap-MatchRows 20 21
Swap-MatchRows 36 37
Swap-MatchRows 48 49
Swap-MatchRows 52 53
Swap-MatchRows 64 65

function Add-MatchRow($row, $indice, $dataPartida, $home, $homeGols, $away, $awayGols, `
    $homeOpenOdds, $homeOpenDH, $homeCloseOdds, $homeCloseDH, `
    $drawOpenOdds, $drawOpenDH, $drawCloseOdds, $drawCloseDH, `
    $awayOpenOdds, $awayOpenDH, $awayCloseOdds, $awayCloseDH, $url) {

    # Copy the formatting (styles/number formats) from the previous row so
    # that Indice (A) and data_partida (E) keep the same styling used
    # throughout the table.
    $ws.Range("A" + ($row - 1) + ":V" + ($row - 1)).Copy()
    $ws.Range("A$row" + ":V$row").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $indice
    $ws.Cells.Item($row, 2).Value = "switzerland"
    $ws.Cells.Item($row, 3).Value = "super-league"
    $ws.Cells.Item($row, 4).Value = "2023-2024"
    $ws.Cells.Item($row, 5).Value = $dataPartida
    $ws.Cells.Item($row, 6).Value = $home
    $ws.Cells.Item($row, 7).Value = $homeGols
    $ws.Cells.Item($row, 8).Value = $away
    $ws.Cells.Item($row, 9).Value = $awayGols
    $ws.Cells.Item($row, 10).Value = $homeOpenOdds
    $ws.Cells.Item($row, 11).Value = $homeOpenDH
    $ws.Cells.Item($row, 12).Value = $homeCloseOdds
    $ws.Cells.Item($row, 13).Value = $homeCloseDH
    $ws.Cells.Item($row, 14).Value = $drawOpenOdds
    $ws.Cells.Item($row, 15).Value = $drawOpenDH
    $ws.Cells.Item($row, 16).Value = $drawCloseOdds
    $ws.Cells.Item($row, 17).Value = $drawCloseDH
    $ws.Cells.Item($row, 18).Value = $awayOpenOdds
    $ws.Cells.Item($row, 19).Value = $awayOpenDH
    $ws.Cells.Item($row, 20).Value = $awayCloseOdds
    $ws.Cells.Item($row, 21).Value = $awayCloseDH
    $ws.Cells.Item($row, 22).Value = $url
}

Add-MatchRow 66 65 45227.75 "St. Gallen" 3 "Grasshoppers" 1 `
    1.39 "22/10/2023 16:42" 1.62 "28/10/2023 17:36" `
    5.14 "22/10/2023 16:42" 4.56 "28/10/2023 17:59" `
    6.33 "22/10/2023 16:42" 5.02 "28/10/2023 17:59" `
    "https://www.betexplorer.com/football/switzerland/super-league/st-gallen-grasshoppers/fNVV8zY3/"

Add-MatchRow 67 66 45227.75 "Yverdon" 1 "Winterthur" 1 `
    2.65 "22/10/2023 16:42" 2.47 "28/10/2023 17:59" `
    3.64 "22/10/2023 16:42" 3.73 "28/10/2023 17:57" `
    2.44 "22/10/2023 16:42" 2.78 "28/10/2023 17:59" `
    "https://www.betexplorer.com/football/switzerland/super-league/yverdon-winterthur/8QZZ7fmA/"

Add-MatchRow 68 67 45227.85416666666 "Zurich" 1 "Lausanne Ouchy" 1 `
    1.49 "21/10/2023 20:43" 1.51 "28/10/2023 20:17" `
    4.67 "21/10/2023 20:43" 4.53 "28/10/2023 20:22" `
    6.15 "21/10/2023 20:43" 6.41 "28/10/2023 20:26" `
    "https://www.betexplorer.com/football/switzerland/super-league/zurich-lausanne-ouchy/2mZw7E2G/"

Add-MatchRow 69 68 45228.59375 "Lausanne" 3 "Basel" 0 `
    1.93 "22/10/2023 15:43" 2.09 "29/10/2023 14:14" `
    3.9 "22/10/2023 15:43" 3.78 "29/10/2023 14:10" `
    3.68 "22/10/2023 15:43" 3.44 "29/10/2023 14:14" `
    "https://www.betexplorer.com/football/switzerland/super-league/lausanne-basel/phlunDQj/"

Add-MatchRow 70 69 45228.6875 "Servette" 4 "Luzern" 2 `
    1.99 "22/10/2023 20:15" 1.93 "29/10/2023 16:21" `
    3.86 "22/10/2023 20:15" 3.87 "29/10/2023 16:21" `
    3.52 "22/10/2023 20:15" 3.9 "29/10/2023 16:21" `
    "https://www.betexplorer.com/football/switzerland/super-league/servette-luzern/CbtR9GIc/"

Add-MatchRow 71 70 45228.6875 "Lugano" 1 "Young Boys" 1 `
    2.95 "22/10/2023 20:15" 2.73 "29/10/2023 16:28" `
    3.68 "22/10/2023 20:15" 3.66 "29/10/2023 16:28" `
    2.22 "22/10/2023 20:15" 2.53 "29/10/2023 16:29" `
    "https://www.betexplorer.com/football/switzerland/super-league/lugano-young-boys/Y7sNAd3i/"
